$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "DHTC"
$ws.Range("B3").Value = "DHTC - Đơn hàng thành công"
